$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 8
$ws.Range("C2").Value = 0.5023643612451429
$ws.Range("D2").Value = 0.131400507317854
$ws.Range("E2").Value = "norm_coldread_gaze_wpm_median"
$ws.Range("F2").Value = ""

# Row 3
$ws.Range("B3").Value = 8
$ws.Range("C3").Value = 0.5017701488129156
$ws.Range("D3").Value = 0.08666302956686936
$ws.Range("E3").Value = "norm_coldread_gaze_wpm_median"
$ws.Range("F3").Value = ""

# Row 4
$ws.Range("B4").Value = 9
$ws.Range("C4").Value = 0.4643702258132781
$ws.Range("D4").Value = 0.1050034032094902
$ws.Range("E4").Value = "norm_coldread_gaze_wpm_median"
$ws.Range("F4").Value = ""
